$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update historical financial data for rows 2-6 (years 2014/12 - 2018/12)
# Values below correct prior erroneous figures reported in the IFRS list

# Row 2
$ws.Range("D2").Value = 1130
$ws.Range("E2").Value = -142
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = -340
$ws.Range("H2").Value = -314
$ws.Range("I2").Value = -323
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 5663
$ws.Range("L2").Value = 5297
$ws.Range("M2").Value = 366
$ws.Range("N2").Value = 308
$ws.Range("O2").Value = 58
$ws.Range("P2").Value = 167
$ws.Range("Q2").Value = -13
$ws.Range("R2").Value = -35
$ws.Range("S2").Value = -8
$ws.Range("T2").Value = 60
$ws.Range("U2").Value = -72
$ws.Range("V2").Value = 754
$ws.Range("W2").Value = -12.52
$ws.Range("X2").Value = -27.77
$ws.Range("Y2").Value = -56.62
$ws.Range("Z2").Value = -5.71
$ws.Range("AA2").Value = 1447.45
$ws.Range("AB2").Value = 213.32
$ws.Range("AC2").Value = -990
$ws.Range("AD2").Value = -1.44
$ws.Range("AE2").Value = 1064
$ws.Range("AF2").Value = 1.34
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 32639841

# Row 3
$ws.Range("D3").Value = 1969
$ws.Range("E3").Value = 168
$ws.Range("F3").Value = 168
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = -175
$ws.Range("I3").Value = -169
$ws.Range("J3").Value = -6
$ws.Range("K3").Value = 6433
$ws.Range("L3").Value = 6182
$ws.Range("M3").Value = 251
$ws.Range("N3").Value = 216
$ws.Range("O3").Value = 36
$ws.Range("P3").Value = 206
$ws.Range("Q3").Value = 351
$ws.Range("R3").Value = -49
$ws.Range("S3").Value = 18
$ws.Range("T3").Value = 61
$ws.Range("U3").Value = 290
$ws.Range("V3").Value = 655
$ws.Range("W3").Value = 8.51
$ws.Range("X3").Value = -8.9
$ws.Range("Y3").Value = -64.5
$ws.Range("Z3").Value = -2.9
$ws.Range("AA3").Value = 2458.97
$ws.Range("AB3").Value = 120.27
$ws.Range("AC3").Value = -509
$ws.Range("AD3").Value = -2.81
$ws.Range("AE3").Value = 588
$ws.Range("AF3").Value = 2.43
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 40409849

# Row 4
$ws.Range("D4").Value = 1848
$ws.Range("E4").Value = 275
$ws.Range("F4").Value = 275
$ws.Range("G4").Value = 216
$ws.Range("H4").Value = 130
$ws.Range("I4").Value = 59
$ws.Range("J4").Value = 71
$ws.Range("K4").Value = 6870
$ws.Range("L4").Value = 6334
$ws.Range("M4").Value = 536
$ws.Range("N4").Value = 427
$ws.Range("O4").Value = 109
$ws.Range("P4").Value = 250
$ws.Range("Q4").Value = -48
$ws.Range("R4").Value = 49
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 17
$ws.Range("U4").Value = -65
$ws.Range("V4").Value = 517
$ws.Range("W4").Value = 14.85
$ws.Range("X4").Value = 7.04
$ws.Range("Y4").Value = 18.45
$ws.Range("Z4").Value = 1.96
$ws.Range("AA4").Value = 1181.79
$ws.Range("AB4").Value = 149.8
$ws.Range("AC4").Value = 122
$ws.Range("AD4").Value = 21.92
$ws.Range("AE4").Value = 916
$ws.Range("AF4").Value = 2.92
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 49189958

# Row 5
$ws.Range("D5").Value = 2763
$ws.Range("E5").Value = 297
$ws.Range("F5").Value = 297
$ws.Range("G5").Value = 185
$ws.Range("H5").Value = 146
$ws.Range("I5").Value = 84
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 7469
$ws.Range("L5").Value = 6820
$ws.Range("M5").Value = 649
$ws.Range("N5").Value = 480
$ws.Range("O5").Value = 168
$ws.Range("P5").Value = 250
$ws.Range("Q5").Value = 403
$ws.Range("R5").Value = 127
$ws.Range("S5").Value = 33
$ws.Range("T5").Value = 15
$ws.Range("U5").Value = 388
$ws.Range("V5").Value = 552
$ws.Range("W5").Value = 10.76
$ws.Range("X5").Value = 5.29
$ws.Range("Y5").Value = 18.42
$ws.Range("Z5").Value = 2.04
$ws.Range("AA5").Value = 1051.26
$ws.Range("AB5").Value = 180.41
$ws.Range("AC5").Value = 170
$ws.Range("AD5").Value = 15.72
$ws.Range("AE5").Value = 1030
$ws.Range("AF5").Value = 2.59
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 49189958

# Row 6
$ws.Range("D6").Value = 2215
$ws.Range("E6").Value = 137
$ws.Range("F6").Value = 137
$ws.Range("G6").Value = 136
$ws.Range("H6").Value = 143
$ws.Range("I6").Value = 96
$ws.Range("K6").Value = 7697
$ws.Range("L6").Value = 6853
$ws.Range("M6").Value = 845
$ws.Range("N6").Value = 628
$ws.Range("P6").Value = 250
$ws.Range("Q6").Value = -238
$ws.Range("R6").Value = 10
$ws.Range("S6").Value = 149
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = -242
$ws.Range("V6").Value = 591
$ws.Range("W6").Value = 6.16
$ws.Range("X6").Value = 6.47
$ws.Range("Y6").Value = 17.39
$ws.Range("Z6").Value = 1.89
$ws.Range("AA6").Value = 811.2
$ws.Range("AB6").Value = 258.63
$ws.Range("AC6").Value = 196
$ws.Range("AD6").Value = 8.24
$ws.Range("AE6").Value = 1289
$ws.Range("AF6").Value = 1.25
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 49189958

# Rows 7-9 correspond to analyst estimate years (2019/12(E), 2020/12(E), 2021/12(E))
# which are no longer populated in this data pull - clear their financial figures
$ws.Range("D7:AJ9").ClearContents()

Write-Output "IFRS list corrected: historical rows updated, estimate rows cleared"
